$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new data rows (6, 7, 8) to the NIFTY options sheet, following the
# exact same layout/formatting as the existing data rows (2-5).
#
# Strategy: first tile-copy the full formatting of row 5 onto rows 6-8 in one
# shot (this both creates the cells and gives them the right style indexes,
# without generating any new/duplicate style entries). Then fill in the
# actual values for every column. Columns that Excel's COM layer would
# otherwise "smart convert" (date-looking text in column A, percent-looking
# text in column E) are temporarily forced to Text format before the value is
# written, then the row-5 number format is re-applied so the final style
# exactly matches the target (General format, like the rest of the sheet).
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

# 1) Clone formatting of row 5 across the three new rows.
$ws.Range("A5:AE5").Copy()
$ws.Range("A6:AE8").PasteSpecial($xlPasteFormats)

# 2) Guard the columns that hold text which looks like a date / percentage
#    so the COM layer stores them as literal text instead of re-interpreting
#    them as a date serial number / percentage number.
$ws.Range("A6:A8").NumberFormat = "@"
$ws.Range("E6:E8").NumberFormat = "@"

# 3) Row 6 - 2026-02-09
$ws.Range("A6").Value = "2026-02-09"
$ws.Range("B6").Value = "10:00:13"
$ws.Range("C6").Value = "AVOID"
$ws.Range("D6").Value = "AVOID"
$ws.Range("E6").Value = "100%"
$ws.Range("F6").Value = "TRADEABLE"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 25814.7
$ws.Range("I6").Value = 12.15
$ws.Range("J6").Value = -0.1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 43.1
$ws.Range("M6").Value = "UNKNOWN"
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = "UNKNOWN"
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = "NONE"
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = "HARD VETO: CPR TRENDING DAY: Price 25814.70 above TC 25661.78 - BULLISH TRENDING DAY likely"
$ws.Range("AD6").Value = "CPR TRENDING DAY: Price 25814.70 above TC 25661.78 - BULLISH TRENDING DAY likely"
$ws.Range("AE6").Value = "Yes"

# 4) Row 7 - 2026-02-10
$ws.Range("A7").Value = "2026-02-10"
$ws.Range("B7").Value = "10:00:10"
$ws.Range("C7").Value = "AVOID"
$ws.Range("D7").Value = "AVOID"
$ws.Range("E7").Value = "100%"
$ws.Range("F7").Value = "TRADEABLE"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 25949.95
$ws.Range("I7").Value = 11.95
$ws.Range("J7").Value = -0.22
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 37.5
$ws.Range("M7").Value = "UNKNOWN"
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = "UNKNOWN"
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = "NONE"
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0
$ws.Range("AC7").Value = "HARD VETO: CPR TRENDING DAY: Price 25949.95 above TC 25862.06 - BULLISH TRENDING DAY likely"
$ws.Range("AD7").Value = "CPR TRENDING DAY: Price 25949.95 above TC 25862.06 - BULLISH TRENDING DAY likely"
$ws.Range("AE7").Value = "Yes"

# 5) Row 8 - 2026-02-11
$ws.Range("A8").Value = "2026-02-11"
$ws.Range("B8").Value = "10:00:08"
$ws.Range("C8").Value = "AVOID"
$ws.Range("D8").Value = "AVOID"
$ws.Range("E8").Value = "100%"
$ws.Range("F8").Value = "TRADEABLE"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 25943.05
$ws.Range("I8").Value = 11.54
$ws.Range("J8").Value = -0.4
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 29.3
$ws.Range("M8").Value = "UNKNOWN"
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = "UNKNOWN"
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = "NONE"
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = "HARD VETO: CPR TRENDING DAY: Price 25943.05 above TC 25933.42 - BULLISH TRENDING DAY likely"
$ws.Range("AD8").Value = "CPR TRENDING DAY: Price 25943.05 above TC 25933.42 - BULLISH TRENDING DAY likely"
$ws.Range("AE8").Value = "Yes"

# 6) Re-apply row 5's number format to columns A and E so their final style
#    matches the rest of the sheet exactly (General format, same style index
#    as every other text cell in the row) instead of the temporary "@" guard.
$ws.Range("A5").Copy()
$ws.Range("A6:A8").PasteSpecial($xlPasteFormats)
$ws.Range("E5").Copy()
$ws.Range("E6:E8").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
